$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '70.646.82'
$ws.Range('E2').Value = '  +0.73%  '

# Row 3
$ws.Range('D3').Value = '3.583.17'
$ws.Range('E3').Value = '  +0.07%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.12%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '585.77'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.49%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '186.16'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.18%  '

# Row 7
$ws.Range('D7').Value = '3.570.59'
$ws.Range('E7').Value = '  -0.11%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.621'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.33%  '

# Row 9
$ws.Range('E9').Value = '  +0.03%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.216'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +17.77%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.650'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.04%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '54.20'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.64%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000323'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +5.68%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.54'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.06%  '

# Row 15
$ws.Range('D15').Value = '4.148.48'
$ws.Range('E15').Value = '  -0.04%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '19.56'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.55%  '

# Row 17
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.608.50'

# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '70.607.96'
$ws.Range('E18').Value = '  +0.90%  '

# Row 19
$ws.Range('E19').Value = '  -1.42%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '567.76'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +15.24%  '

# Row 21
$ws.Range('E21').Value = '  -0.26%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.01'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.91%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '17.68'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -7.38%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.65'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +6.04%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.90'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.93%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '95.69'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.23%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.50'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.87%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.94'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.65%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.13'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.21%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '32.27'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.86%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.30'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -6.07%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '12.44'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +2.68%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '64.82'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.68%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.114'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.74%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.33'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.29%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '562.79'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.93%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.418'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.14%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '37.74'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.95%  '

# Row 39
$ws.Range('E39').Value = '  +0.15%  '

# Row 40
$ws.Range('D40').Value = '0.0₃0799'
$ws.Range('E40').Value = '  +0.83%  '

# Row 41
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.14'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -1.92%  '

# Row 42
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '3.376.55'
$ws.Range('E42').Value = '  +5.41%  '

# Row 43
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.134'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.52%  '

# Row 44
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.37'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.60%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.57'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.80%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0445'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.53%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.96'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.46%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.49'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.55%  '

# Row 49
$ws.Range('E49').Value = '  +0.60%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.999'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.12%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.44'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -10.87%  '
